# Add a new record (row 16) to the hardrock2015 sheet representing an
# effort that is missing the StartSplit time but has times recorded at
# later splits ("Distance from Start" and "Ridgeline In").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A16").Value = "Missing"
$ws.Range("B16").Value = "StartSplit"
$ws.Range("C16").Value = "M"
$ws.Range("D16").Value = "NY"
$ws.Range("H16").Value = 11
$ws.Range("J16").Value = 0.14583333333333334
$ws.Range("K16").Value = 0.14930555555555555

# Match the time-of-day number format used by the other split-time columns.
$ws.Range("J16:K16").NumberFormat = "h:mm"

# Move the active selection to L16, matching the saved workbook state.
$ws.Range("L16").Select() | Out-Null
